$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.522.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.985.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.23%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.460.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.984.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.561.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.171"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.97%  "
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.60%  "
$ws.Range("E43").Value = "  +3.03%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.270"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.284.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.536"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +19.67%  "
